$wb = $excel.ActiveWorkbook

# ALC row 28 (G28=27772)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 599.9
$ws.Range("I28").Value = 671.4286
$ws.Range("J28").Value = 433
$ws.Range("K28").Value = 671.4286
$ws.Range("L28").Value = 433
$ws.Range("M28").Value = -186.4286
$ws.Range("N28").Value = -1403

# ALC row 61 (G61=4604)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H61").Value = 98.333336
$ws.Range("I61").Value = 98.333336
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 295.000008
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -123.000008
$ws.Range("N61").ClearContents()

# ALC row 106 (G106=19903)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 2235
$ws.Range("I106").Value = 2187.8572
$ws.Range("K106").Value = 2187.8572
$ws.Range("M106").Value = -1556.8572

# ARM row 30 (G30=2712)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("M30").ClearContents()
$ws.Range("N30").ClearContents()

# ARM row 32 (G32=44147)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5637.82
$ws.Range("I32").Value = 4425.643
$ws.Range("J32").Value = 12001.75
$ws.Range("K32").Value = 4425.643
$ws.Range("L32").Value = 12001.75
$ws.Range("M32").Value = -4138.643
$ws.Range("N32").Value = -12575.75

# ARM row 34 (G34=2753)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()

# ARM row 35 (G35=2473)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()

# ARM row 42 (G42=2765)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H42").Value = 7565.2856
$ws.Range("J42").Value = 7565.2856
$ws.Range("L42").Value = 7565.2856
$ws.Range("N42").Value = -8537.285599999999

# ARM row 45 (G45=27714)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1550.75
$ws.Range("I45").Value = 1606.5264
$ws.Range("K45").Value = 1606.5264
$ws.Range("M45").Value = -1229.5264

# ARM row 61 (G61=43999)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1896.2858
$ws.Range("I61").Value = 1748.3572
$ws.Range("J61").Value = 2488
$ws.Range("K61").Value = 1748.3572
$ws.Range("L61").Value = 2488
$ws.Range("M61").Value = -1536.3572
$ws.Range("N61").Value = -2912

# ARM row 74 (G74=44000)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 388295.03
$ws.Range("I74").Value = 4566.357
$ws.Range("J74").Value = 835978.5
$ws.Range("K74").Value = 4566.357
$ws.Range("L74").Value = 835978.5
$ws.Range("M74").Value = -3692.357
$ws.Range("N74").Value = -837726.5

# ARM row 77 (G77=44000)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 388295.03
$ws.Range("I77").Value = 4566.357
$ws.Range("J77").Value = 835978.5
$ws.Range("K77").Value = 22831.785
$ws.Range("L77").Value = 4179892.5
$ws.Range("M77").Value = -18463.785
$ws.Range("N77").Value = -4188628.5

# ARM row 107 (G107=25645)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H107").Value = 50000
$ws.Range("J107").Value = 50000
$ws.Range("L107").Value = 50000
$ws.Range("N107").Value = -57680

# ARM row 109 (G109=25646)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H109").Value = 30000
$ws.Range("J109").Value = 30000
$ws.Range("L109").Value = 30000
$ws.Range("N109").Value = -32774

# ARM row 110 (G110=27708)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 3140.923
$ws.Range("I110").Value = 3791.625
$ws.Range("J110").Value = 2099.8
$ws.Range("K110").Value = 3791.625
$ws.Range("L110").Value = 2099.8
$ws.Range("M110").Value = -1746.625
$ws.Range("N110").Value = -6189.8

# ARM row 136 (G136=43999)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1896.2858
$ws.Range("I136").Value = 1748.3572
$ws.Range("J136").Value = 2488
$ws.Range("K136").Value = 5245.071599999999
$ws.Range("L136").Value = 7464
$ws.Range("M136").Value = -2695.071599999999
$ws.Range("N136").Value = -12564

# BSM row 30 (G30=1609)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()

# BSM row 61 (G61=2543)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H61").Value = 30106
$ws.Range("J61").Value = 30106
$ws.Range("L61").Value = 30106
$ws.Range("N61").Value = -30732

# CRP row 8 (G8=1894)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 1045
$ws.Range("I8").Value = 600
$ws.Range("J8").Value = 1490
$ws.Range("K8").Value = 600
$ws.Range("L8").Value = 1490
$ws.Range("M8").Value = -460
$ws.Range("N8").Value = -1770

# CRP row 16 (G16=27691)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2011.091
$ws.Range("I16").Value = 1980.2222
$ws.Range("K16").Value = 1980.2222
$ws.Range("M16").Value = -1693.2222

# CRP row 22 (G22=5367)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 430.75
$ws.Range("I22").Value = 336
$ws.Range("J22").Value = 588.6667
$ws.Range("K22").Value = 336
$ws.Range("L22").Value = 588.6667
$ws.Range("M22").Value = 14
$ws.Range("N22").Value = -1288.6667

# CRP row 32 (G32=2246)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value = 1672.8572
$ws.Range("I32").Value = 1618.3334
$ws.Range("J32").Value = 2000
$ws.Range("K32").Value = 1618.3334
$ws.Range("L32").Value = 2000
$ws.Range("M32").Value = -1302.3334
$ws.Range("N32").Value = -2632

# CRP row 35 (G35=1627)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()

# CRP row 113 (G113=27691)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 2011.091
$ws.Range("I113").Value = 1980.2222
$ws.Range("K113").Value = 1980.2222
$ws.Range("M113").Value = 189.7778000000001

# CRP row 132 (G132=44019)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2054.6
$ws.Range("I132").Value = 1844.1364
$ws.Range("J132").Value = 2633.375
$ws.Range("K132").Value = 5532.4092
$ws.Range("L132").Value = 7900.125
$ws.Range("M132").Value = -3002.4092
$ws.Range("N132").Value = -12960.125

# CUL row 68 (G68=12895)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1244.9048
$ws.Range("I68").Value = 800
$ws.Range("J68").Value = 1383.9375
$ws.Range("K68").Value = 2400
$ws.Range("L68").Value = 4151.8125
$ws.Range("M68").Value = -1589
$ws.Range("N68").Value = -5773.8125

# CUL row 71 (G71=12895)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 1244.9048
$ws.Range("I71").Value = 800
$ws.Range("J71").Value = 1383.9375
$ws.Range("K71").Value = 7200
$ws.Range("L71").Value = 12455.4375
$ws.Range("M71").Value = -3144
$ws.Range("N71").Value = -20567.4375

# CUL row 113 (G113=27843)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1165.2941
$ws.Range("I113").Value = 439.4
$ws.Range("J113").Value = 1467.75
$ws.Range("K113").Value = 1318.2
$ws.Range("L113").Value = 4403.25
$ws.Range("M113").Value = 851.8000000000002
$ws.Range("N113").Value = -8743.25

# GSM row 29 (G29=4209)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 2000
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 2000
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 2000
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -2580

# GSM row 102 (G102=36169)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1229.2273
$ws.Range("I102").Value = 1082.8667
$ws.Range("J102").Value = 1542.8572
$ws.Range("K102").Value = 1082.8667
$ws.Range("L102").Value = 1542.8572
$ws.Range("M102").Value = 539.1333
$ws.Range("N102").Value = -4786.8572

# GSM row 132 (G132=44008)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1835.7959
$ws.Range("I132").Value = 1297.7097
$ws.Range("J132").Value = 2762.5
$ws.Range("K132").Value = 3893.1291
$ws.Range("L132").Value = 8287.5
$ws.Range("M132").Value = -1363.1291
$ws.Range("N132").Value = -13347.5

# LTW row 16 (G16=5289)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 282.1579
$ws.Range("I16").Value = 196
$ws.Range("J16").Value = 523.4
$ws.Range("K16").Value = 196
$ws.Range("L16").Value = 523.4
$ws.Range("M16").Value = -26
$ws.Range("N16").Value = -863.4

# LTW row 22 (G22=5277)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4167854
$ws.Range("I22").Value = 8333708
$ws.Range("J22").Value = 2000
$ws.Range("K22").Value = 8333708
$ws.Range("L22").Value = 2000
$ws.Range("M22").Value = -8333413
$ws.Range("N22").Value = -2590

# LTW row 27 (G27=5277)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 4167854
$ws.Range("I27").Value = 8333708
$ws.Range("J27").Value = 2000
$ws.Range("K27").Value = 8333708
$ws.Range("L27").Value = 2000
$ws.Range("M27").Value = -8333601
$ws.Range("N27").Value = -2214

# LTW row 38 (G38=2767)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 10541.5
$ws.Range("I38").Value = 7000
$ws.Range("J38").Value = 11722
$ws.Range("K38").Value = 7000
$ws.Range("L38").Value = 11722
$ws.Range("M38").Value = -6590
$ws.Range("N38").Value = -12542

# LTW row 61 (G61=27740)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4151.909
$ws.Range("I61").Value = 3884
$ws.Range("K61").Value = 3884
$ws.Range("M61").Value = -3682

# LTW row 93 (G93=19993)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1686
$ws.Range("I93").Value = 1641.091
$ws.Range("J93").Value = 1756.5714
$ws.Range("K93").Value = 1641.091
$ws.Range("L93").Value = 1756.5714
$ws.Range("M93").Value = -393.0909999999999
$ws.Range("N93").Value = -4252.5714

# LTW row 113 (G113=27740)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 4151.909
$ws.Range("I113").Value = 3884
$ws.Range("K113").Value = 3884
$ws.Range("M113").Value = -1714

# LTW row 132 (G132=44058)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 16676209
$ws.Range("I132").Value = 29425660
$ws.Range("J132").Value = 3849.4614
$ws.Range("K132").Value = 88276980
$ws.Range("L132").Value = 11548.3842
$ws.Range("M132").Value = -88274450
$ws.Range("N132").Value = -16608.3842

# LTW row 136 (G136=44060)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 6408.207
$ws.Range("I136").Value = 9557.267
$ws.Range("J136").Value = 3034.2144
$ws.Range("K136").Value = 28671.801
$ws.Range("L136").Value = 9102.643199999999
$ws.Range("M136").Value = -26121.801
$ws.Range("N136").Value = -14202.6432

# LTW row 140 (G140=42503)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# WVR row 107 (G107=27746)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1005.0323
$ws.Range("I107").Value = 782.6
$ws.Range("J107").Value = 1409.4546
$ws.Range("K107").Value = 2347.8
$ws.Range("L107").Value = 4228.3638
$ws.Range("M107").Value = -427.8000000000002
$ws.Range("N107").Value = -8068.3638

# WVR row 113 (G113=27752)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 422.15
$ws.Range("I113").Value = 313.8
$ws.Range("J113").Value = 747.2
$ws.Range("K113").Value = 941.4000000000001
$ws.Range("L113").Value = 2241.6
$ws.Range("M113").Value = 1228.6
$ws.Range("N113").Value = -6581.6

# WVR row 132 (G132=44029)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1388.8182
$ws.Range("I132").Value = 787.56757
$ws.Range("J132").Value = 2624.7222
$ws.Range("K132").Value = 2362.70271
$ws.Range("L132").Value = 7874.1666
$ws.Range("M132").Value = 167.29729
$ws.Range("N132").Value = -12934.1666

# WVR row 136 (G136=44031)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4123.122
$ws.Range("I136").Value = 1366.7142
$ws.Range("J136").Value = 10060
$ws.Range("K136").Value = 4100.142599999999
$ws.Range("L136").Value = 30180
$ws.Range("M136").Value = -1550.142599999999
$ws.Range("N136").Value = -35280
